$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (existing A16-E16, F16, G16, K16, L16). Fill in the missing H/I/J/M/N values.
$ws.Range("H16").Value = 26228
$ws.Range("I16").Value = 0.26866099999999998
$ws.Range("J16").Value = 127027
$ws.Range("M16").Value = 1490
$ws.Range("N16").Value = 3.655243

# Row 17 (existing A17, F17, K17, L17). Add B-E, G-J, M, N.
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 30303
$ws.Range("H17").Value = 5560
$ws.Range("I17").Value = 0.84496599999999999
$ws.Range("J17").Value = 314529
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 9.0506729999999997

# Row 18 (existing A18, F18, K18, L18). Add B-E, G-J, M, N.
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 1
$ws.Range("G18").Value = 30303
$ws.Range("H18").Value = 5560
$ws.Range("I18").Value = 0.84496599999999999
$ws.Range("J18").Value = 307485
$ws.Range("M18").Value = 7044
$ws.Range("N18").Value = 8.8479799999999997

# Row 19 (existing A19, F19, K19, L19). Add B-E, G-J, M, N.
$ws.Range("B19").Value = 4
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 35766
$ws.Range("H19").Value = 97
$ws.Range("I19").Value = 0.99729500000000004
$ws.Range("J19").Value = 363696
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 10.46547

# Row 20 (existing A20, F20, K20, L20). Add B-E, G-J, M, N.
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 1
$ws.Range("G20").Value = 35766
$ws.Range("H20").Value = 97
$ws.Range("I20").Value = 0.99729500000000004
$ws.Range("J20").Value = 356652
$ws.Range("M20").Value = 7044
$ws.Range("N20").Value = 10.262776000000001

# Update the selected cell (per the diff, the active cell/selection moves to H21).
$ws.Range("H21").Select()
